$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 209.41667
$ws.Range("I11").Value = 209.41667
$ws.Range("K11").Value = 209.41667
$ws.Range("M11").Value = -69.41667000000001
$ws.Range("H15").Value = 1494.5
$ws.Range("I15").Value = 1494.5
$ws.Range("K15").Value = 4483.5
$ws.Range("M15").Value = -4314.5
$ws.Range("H87").Value = 19833.334
$ws.Range("J87").Value = 19833.334
$ws.Range("L87").Value = 19833.334
$ws.Range("N87").Value = -22329.334
$ws.Range("H90").Value = 19833.334
$ws.Range("J90").Value = 19833.334
$ws.Range("L90").Value = 59500.00199999999
$ws.Range("N90").Value = -71980.00199999999
$ws.Range("H138").Value = 4671.6177
$ws.Range("I138").Value = 3995.3333
$ws.Range("K138").Value = 11985.9999
$ws.Range("M138").Value = -6845.999899999999

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 19593.895
$ws.Range("I2").Value = 23452.8
$ws.Range("J2").Value = 5123
$ws.Range("K2").Value = 23452.8
$ws.Range("L2").Value = 5123
$ws.Range("M2").Value = -23339.8
$ws.Range("N2").Value = -5349
$ws.Range("H61").Value = 4236.0967
$ws.Range("I61").Value = 3096.125
$ws.Range("J61").Value = 8144.5713
$ws.Range("K61").Value = 3096.125
$ws.Range("L61").Value = 8144.5713
$ws.Range("M61").Value = -2884.125
$ws.Range("N61").Value = -8568.5713
$ws.Range("H110").Value = 3111.5557
$ws.Range("I110").Value = 2388.318
$ws.Range("K110").Value = 2388.318
$ws.Range("M110").Value = -343.3180000000002
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("H116").Value = 19593.895
$ws.Range("I116").Value = 23452.8
$ws.Range("J116").Value = 5123
$ws.Range("K116").Value = 23452.8
$ws.Range("L116").Value = 5123
$ws.Range("M116").Value = -21158.8
$ws.Range("N116").Value = -9711
$ws.Range("H117").Value = 100000
$ws.Range("J117").Value = 100000
$ws.Range("L117").Value = 100000
$ws.Range("N117").Value = -109178
$ws.Range("H118").Value = 20454.545
$ws.Range("J118").Value = 20454.545
$ws.Range("L118").Value = 20454.545
$ws.Range("N118").Value = -23768.545
$ws.Range("H122").Value = 3235.1226
$ws.Range("I122").Value = 2907.1396
$ws.Range("J122").Value = 5585.6665
$ws.Range("K122").Value = 8721.418799999999
$ws.Range("L122").Value = 16756.9995
$ws.Range("M122").Value = -6271.418799999999
$ws.Range("N122").Value = -21656.9995
$ws.Range("H123").Value = 30000
$ws.Range("J123").Value = 30000
$ws.Range("L123").Value = 30000
$ws.Range("N123").Value = -39800
$ws.Range("H125").Value = 40000
$ws.Range("J125").Value = 40000
$ws.Range("L125").Value = 40000
$ws.Range("N125").Value = -49840
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("H130").Value = 776000
$ws.Range("J130").Value = 776000
$ws.Range("L130").Value = 776000
$ws.Range("N130").Value = -786040
$ws.Range("H131").Value = 42724.273
$ws.Range("J131").Value = 42724.273
$ws.Range("L131").Value = 42724.273
$ws.Range("N131").Value = -52804.273
$ws.Range("H132").Value = 3084.2407
$ws.Range("I132").Value = 1719.4884
$ws.Range("J132").Value = 8419.182000000001
$ws.Range("K132").Value = 5158.4652
$ws.Range("L132").Value = 25257.546
$ws.Range("M132").Value = -2628.4652
$ws.Range("N132").Value = -30317.546
$ws.Range("H135").Value = 50000
$ws.Range("J135").Value = 50000
$ws.Range("L135").Value = 50000
$ws.Range("N135").Value = -60140
$ws.Range("H136").Value = 4236.0967
$ws.Range("I136").Value = 3096.125
$ws.Range("J136").Value = 8144.5713
$ws.Range("K136").Value = 9288.375
$ws.Range("L136").Value = 24433.7139
$ws.Range("M136").Value = -6738.375
$ws.Range("N136").Value = -29533.7139
$ws.Range("N111").ClearContents()
$ws.Range("N115").ClearContents()
$ws.Range("N127").ClearContents()
$ws.Range("N129").ClearContents()

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 19593.895
$ws.Range("I3").Value = 23452.8
$ws.Range("J3").Value = 5123
$ws.Range("K3").Value = 23452.8
$ws.Range("L3").Value = 5123
$ws.Range("M3").Value = -23338.8
$ws.Range("N3").Value = -5351
$ws.Range("H64").Value = 1927
$ws.Range("J64").Value = 2372.25
$ws.Range("L64").Value = 2372.25
$ws.Range("N64").Value = -2822.25
$ws.Range("H67").Value = 1927
$ws.Range("J67").Value = 2372.25
$ws.Range("L67").Value = 2372.25
$ws.Range("N67").Value = -3932.25
$ws.Range("H86").Value = 49620.285
$ws.Range("I86").Value = 73055.71000000001
$ws.Range("K86").Value = 73055.71000000001
$ws.Range("M86").Value = -71932.71000000001
$ws.Range("H89").Value = 49620.285
$ws.Range("I89").Value = 73055.71000000001
$ws.Range("K89").Value = 365278.55
$ws.Range("M89").Value = -359662.55
$ws.Range("H134").Value = 5515.4287
$ws.Range("I134").Value = 2739.9412
$ws.Range("K134").Value = 8219.8236
$ws.Range("M134").Value = -5684.8236

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 418573.84
$ws.Range("I22").Value = 2642.375
$ws.Range("J22").Value = 1250436.8
$ws.Range("K22").Value = 2642.375
$ws.Range("L22").Value = 1250436.8
$ws.Range("M22").Value = -2292.375
$ws.Range("N22").Value = -1251136.8
$ws.Range("H38").Value = 11657.6
$ws.Range("I38").Value = 8322
$ws.Range("J38").Value = 25000
$ws.Range("K38").Value = 8322
$ws.Range("L38").Value = 25000
$ws.Range("M38").Value = -7945
$ws.Range("N38").Value = -25754
$ws.Range("H42").Value = 17750
$ws.Range("I42").Value = 5000
$ws.Range("J42").Value = 22000
$ws.Range("K42").Value = 5000
$ws.Range("L42").Value = 22000
$ws.Range("M42").Value = -4407
$ws.Range("N42").Value = -23186
$ws.Range("H44").Value = 7500
$ws.Range("I44").Value = 5000
$ws.Range("J44").Value = 10000
$ws.Range("K44").Value = 5000
$ws.Range("L44").Value = 10000
$ws.Range("M44").Value = -4558
$ws.Range("N44").Value = -10884
$ws.Range("H45").Value = 30000
$ws.Range("J45").Value = 30000
$ws.Range("L45").Value = 30000
$ws.Range("N45").Value = -31186
$ws.Range("H46").Value = 11657.6
$ws.Range("I46").Value = 8322
$ws.Range("J46").Value = 25000
$ws.Range("K46").Value = 8322
$ws.Range("L46").Value = 25000
$ws.Range("M46").Value = -8111
$ws.Range("N46").Value = -25422
$ws.Range("H47").Value = 27825
$ws.Range("J47").Value = 35000
$ws.Range("L47").Value = 35000
$ws.Range("N47").Value = -36132
$ws.Range("H92").Value = 58497.5
$ws.Range("J92").Value = 58497.5
$ws.Range("L92").Value = 58497.5
$ws.Range("N92").Value = -63489.5

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1480.9231
$ws.Range("J68").Value = 1643.2858
$ws.Range("L68").Value = 4929.857400000001
$ws.Range("N68").Value = -6551.857400000001
$ws.Range("H71").Value = 1480.9231
$ws.Range("J71").Value = 1643.2858
$ws.Range("L71").Value = 14789.5722
$ws.Range("N71").Value = -22901.5722
$ws.Range("H121").Value = 2747.75
$ws.Range("J121").Value = 3854.75
$ws.Range("L121").Value = 11564.25
$ws.Range("N121").Value = -14184.25
$ws.Range("H131").Value = 4038
$ws.Range("I131").Value = 845.6
$ws.Range("K131").Value = 2536.8
$ws.Range("M131").Value = 2503.2

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 26832.166
$ws.Range("J52").Value = 26832.166
$ws.Range("L52").Value = 26832.166
$ws.Range("N52").Value = -27350.166
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("H124").Value = 39769.23
$ws.Range("J124").Value = 39769.23
$ws.Range("L124").Value = 39769.23
$ws.Range("N124").Value = -49589.23
$ws.Range("N93").ClearContents()

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1815.5
$ws.Range("I16").Value = 1874.3334
$ws.Range("J16").Value = 1739.8572
$ws.Range("K16").Value = 1874.3334
$ws.Range("L16").Value = 1739.8572
$ws.Range("M16").Value = -1704.3334
$ws.Range("N16").Value = -2079.8572
$ws.Range("H133").Value = 67959
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 67959
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 67959
$ws.Range("N133").Value = -73019
$ws.Range("H137").Value = 69769.234
$ws.Range("J137").Value = 69769.234
$ws.Range("L137").Value = 69769.234
$ws.Range("N137").Value = -79969.234
$ws.Range("H140").Value = 69997.5
$ws.Range("J140").Value = 69997.5
$ws.Range("L140").Value = 69997.5
$ws.Range("N140").Value = -80357.5
$ws.Range("M133").ClearContents()

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 73409
$ws.Range("J46").Value = 73409
$ws.Range("L46").Value = 73409
$ws.Range("N46").Value = -73871
$ws.Range("H134").Value = 73409
$ws.Range("J134").Value = 73409
$ws.Range("L134").Value = 220227
$ws.Range("N134").Value = -225297
$ws.Range("H135").Value = 71666.664
$ws.Range("J135").Value = 71666.664
$ws.Range("L135").Value = 71666.664
$ws.Range("N135").Value = -81806.664
